$wb = $excel.ActiveWorkbook

# --- 1) Reorder the GL_Acc_Trans_* tabs: Cash moves to sit before Transfer ---
$wsCash = $wb.Worksheets.Item("GL_Acc_Trans_Cash")
$wsTransfer = $wb.Worksheets.Item("GL_Acc_Trans_Transfer")
$wsCash.Move($wsTransfer)

# Re-resolve the handle: Move() leaves old variables pointing at stale slots.
$wsTransfer = $wb.Worksheets.Item("GL_Acc_Trans_Transfer")

# --- 2) Add the new "Share_Transfer" sheet right after GL_Acc_Trans_Transfer (now last) ---
$wsTransfer.Copy($null, $wsTransfer)
$newSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet.Name = "Share_Transfer"

# Drop the Trans_Amount column (old column E) so later columns shift left,
# matching the narrower Share_Transfer layout.
$newSheet.Columns.Item(5).Delete()

# --- 3) Rewrite the header row (row 1) ---
$newSheet.Range("A1").Value = "TestScenario"
$newSheet.Range("B1").Value = "Run"
$newSheet.Range("C1").Value = "pcRegFormName"
$newSheet.Range("D1").Value = "pcRegFormPcName"
$newSheet.Range("E1").Value = "Account_No"
$newSheet.Range("F1").Value = "Account_No1"
$newSheet.Range("G1").Value = ""

# --- 4) Rewrite the data row (row 2) ---
$newSheet.Range("A2").Value = "Share_Transfer"
$newSheet.Range("B2").Value = "Yes"
$newSheet.Range("C2").Value = "qwerty"
$newSheet.Range("D2").Value = "zxcvb"
$newSheet.Range("E2").Value = 86
$newSheet.Range("F2").Value = 87

# --- 5) Update sheet view / selection on the new sheet ---
$newSheet.Range("F11").Select()

# --- 6) Make the new sheet the active tab ---
$newSheet.Activate()
